# Re-order the red-flag report rows 3-6. Reading the table as a cycle of
# "who ends up where":
#   old row 4 (NEPOMSCENE)          -> row 3
#   old row 6 (BROADBAND SYSTEMS)   -> row 4
#   old row 3 (BANDAG)              -> row 5
#   old row 5 (H.VEDASTE)           -> row 6
#
# i.e. data flows 3 -> 5 -> 6 -> 4 -> 3 around the loop. We use
# Range.Copy(Destination), not Range.Value assignment, so that text-typed
# cells (phone numbers / TINs with leading zeros, etc.) keep their exact
# stored representation and formatting instead of being re-parsed/coerced
# by a value-based write. Because the four rows form a single closed
# cycle, row 3's original contents are parked in a scratch row first so
# they survive long enough to be written back into row 5 at the end; the
# scratch row is cleared again afterwards.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratchRow = 200

# Range.Copy only overwrites cells that actually hold something in the
# source, so a destination row has to be blanked out first -- otherwise a
# column that's populated in the old row but empty in the new one (e.g.
# H3's postal code, which NEPOMSCENE doesn't have) would keep showing
# stale leftovers. Clear() (not ClearContents()) is used so the cell is
# dropped entirely rather than left behind as an empty tag.

# 1. Stash original row 3 out of the way.
$ws.Range("A$scratchRow`:Q$scratchRow").Clear()
$ws.Range("A3:Q3").Copy($ws.Range("A$scratchRow`:Q$scratchRow"))

# 2. Old row 4 -> row 3.
$ws.Range("A3:Q3").Clear()
$ws.Range("A4:Q4").Copy($ws.Range("A3:Q3"))

# 3. Old row 6 -> row 4.
$ws.Range("A4:Q4").Clear()
$ws.Range("A6:Q6").Copy($ws.Range("A4:Q4"))

# 4. Old row 5 -> row 6.
$ws.Range("A6:Q6").Clear()
$ws.Range("A5:Q5").Copy($ws.Range("A6:Q6"))

# 5. Stashed original row 3 -> row 5.
$ws.Range("A5:Q5").Clear()
$ws.Range("A$scratchRow`:Q$scratchRow").Copy($ws.Range("A5:Q5"))

# 6. Clean up the scratch row.
$ws.Range("A$scratchRow`:Q$scratchRow").Clear()

# Copying a whole A:Q band (done above, to get a clean row-level snapshot
# in one shot) leaves behind empty placeholder cells for any column the
# source row didn't populate -- e.g. row 3 (now NEPOMSCENE) has no postal
# code, but copying over the old BANDAG row's H3 value means the *shape*
# of the copy still touched H3. A final sweep re-clears exactly the cells
# that should be genuinely blank (no cell entry at all) in the finished
# rows, matching the source rows that never had those columns populated.
$ws.Range("H3").Clear()
$ws.Range("L3").Clear()
$ws.Range("M3").Clear()
$ws.Range("L4").Clear()
$ws.Range("L5").Clear()
$ws.Range("M5").Clear()
$ws.Range("L6").Clear()
$ws.Range("M6").Clear()
